$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5, column B currently holds "Kurek" (duplicate row issue) -> change to "Jurek"
$ws.Range("B5").Value = "Jurek"

# Update the active selection to B15
$ws.Range("B15").Select()
